$d = $word.ActiveDocument

# Replace the date paragraph text "January 2023" with "1/1/23".
$d.Content.Find.Execute("January 2023", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1/1/23", 2)
